$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.186.30'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("E2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.378.55'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.99'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("E5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.78'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.53%  '
$ws.Range("E9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.47'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.23%  '
$ws.Range("E11").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0791'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.79%  '
$ws.Range("E12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.48'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.80'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("E14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.738.25'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("E15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.353.93'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("E16").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.806'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("E17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.170.13'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E18").NumberFormat = "General"
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.35'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.20%  '
$ws.Range("E19").NumberFormat = "General"
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.00'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("E20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0892'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("E21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.16'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.19'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.13%  '
$ws.Range("E23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.34%  '
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.54'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.72%  '
$ws.Range("E28").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.37'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.91%  '
$ws.Range("E29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.00'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.05%  '
$ws.Range("E30").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.11%  '
$ws.Range("E31").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.07'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.31%  '
$ws.Range("E32").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +11.98%  '
$ws.Range("E33").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.92'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.38%  '
$ws.Range("E34").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.19%  '
$ws.Range("E35").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '128.31'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.24%  '
$ws.Range("E36").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.83'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E37").NumberFormat = "General"
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("E38").NumberFormat = "General"
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.85'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.52%  '
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.88%  '
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("E41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.92'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.69%  '
$ws.Range("E42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.931.55'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.28%  '
$ws.Range("E43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0279'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.92%  '
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.77'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.02%  '
$ws.Range("E46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.24'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -7.33%  '
$ws.Range("E47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.596.15'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("E49").NumberFormat = "General"
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.97'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.12%  '
$ws.Range("E50").NumberFormat = "General"
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '71.79'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.49%  '
$ws.Range("E51").NumberFormat = "General"
